# fix(removing title): removing title
# Removes the "CaixaDeTexto 13" textbox ("GOOGLE ADWORDS | OUTUBRO 2022")
# from slide 10 of the presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "CaixaDeTexto 13") {
        $shape.Delete()
    }
}
